# Update Data Sources from LFX
# Re-applies the current (default) table style to every table on every slide.
# This normalizes the table style GUID from
#   {928D4016-75EB-4100-83A0-C58B1C6B50E9}
# to
#   {2148295F-C88A-4D38-9292-040BEF2578C8}
# wherever it is used, matching the source deck's re-export.

$p = $ppt.ActivePresentation

$oldStyleId = "{928D4016-75EB-4100-83A0-C58B1C6B50E9}"
$newStyleId = "{2148295F-C88A-4D38-9292-040BEF2578C8}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
